$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.793.54"
$ws.Range("E2").Value = "  +4.09%  "
$ws.Range("D3").Value = "1.910.30"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.27"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.697"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.80"
$ws.Range("E8").Value = "  +7.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.374"
$ws.Range("E9").Value = "  +4.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.98"
$ws.Range("E10").Value = "  +6.62%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.74"
$ws.Range("E13").Value = "  +8.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.813"
$ws.Range("E14").Value = "  +5.09%  "
$ws.Range("D15").Value = "2.185.41"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D17").Value = "1.919.77"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "36.741.97"
$ws.Range("E18").Value = "  +4.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.37"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("E20").Value = "  +3.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.65"
$ws.Range("E21").Value = "  +6.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "250.64"
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.15"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  -4.37%  "
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.07"
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.60"
$ws.Range("E31").Value = "  +7.24%  "
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.93"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0881"
$ws.Range("E35").Value = "  +19.96%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.21"
$ws.Range("E37").Value = "  +60.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.49"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "104.51"
$ws.Range("E41").Value = "  +7.43%  "
$ws.Range("E42").Value = "  +3.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.83"
$ws.Range("E43").Value = "  +3.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("E44").Value = "  +18.89%  "
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").Value = "1.352.71"
$ws.Range("E46").Value = "  +3.56%  "
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0818"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.81"
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.39"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").Value = "2.087.87"
$ws.Range("E51").Value = "  +1.33%  "
